# session4/Score_iterations.xlsx -- "code returns run len of string"
#
# Appends a new block of log entries (rows 90-97) describing work on:
#   - a follow-up to the previous "n^2 solution" item (row 90 F/G)
#   - "Print list items containing all characters of a given word" (rows 91-93)
#   - a new "Strings" topic header (row 87) followed by:
#       "Reverse words in a given string" (row 94)
#       "Run Length Encoding" (rows 95-97)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New "Strings" topic header -------------------------------------------
$ws.Range("A87").Value = "Strings"
$ws.Range("A87").Font.Bold = $true

# --- finish off the previous (Matrix/array) entry on row 90 ---------------
$ws.Range("F90").Value = "n2 solution"
$ws.Range("G90").Value = "need to find O(n)"

# --- row 91: discussion -----------------------------------------------------
$ws.Range("C91").Value = "discussing above code"
$ws.Range("E91").Value = "x - 11:45"

# --- row 92-93: "Print list items containing all characters of a given word"
$urlPrintList = "https://www.geeksforgeeks.org/print-list-items-containing-all-characters-of-a-given-word/"
$ws.Hyperlinks.Add($ws.Range("A92"), $urlPrintList)
$ws.Range("A92").Style = "Hyperlink"
$ws.Range("B92").Value = "GFG"
$ws.Range("C92").Value = "Print list items containing all characters of a given word"
$ws.Range("D92").Value = "not done, skipped"
$ws.Range("E92").Value = "11:55 - 12:40"

$ws.Range("D93").Value = "frustrated break"

# --- row 94: "Reverse words in a given string" -----------------------------
$urlReverseWords = "https://www.geeksforgeeks.org/reverse-words-in-a-given-string/"
$ws.Hyperlinks.Add($ws.Range("A94"), $urlReverseWords)
$ws.Range("A94").Style = "Hyperlink"
$ws.Range("B94").Value = "GFG"
$ws.Range("C94").Value = "Reverse words in a given string"
$ws.Range("D94").Value = "not done"
$ws.Range("E94").Value = "11:28 - 12:07"

# --- row 95-97: "Run Length Encoding" ---------------------------------------
$urlRunLength = "https://www.geeksforgeeks.org/run-length-encoding/"
$ws.Hyperlinks.Add($ws.Range("A95"), $urlRunLength)
$ws.Range("A95").Style = "Hyperlink"
$ws.Range("B95").Value = "GFG"
$ws.Range("C95").Value = "Run Length Encoding"

$ws.Range("D96").Value = "headache break"

$ws.Range("E97").Value = "1:42 - 1:53"

# --- keep the sheet's view in sync with the new bottom of the log ----------
$ws.Activate() | Out-Null
$ws.Range("E98").Select() | Out-Null
